$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date field text (22/01/2017 ->
#    23/01/2017) on the slide master and every slide layout's date
#    placeholder - mirrors PowerPoint re-caching the auto date field a day
#    later.
# ---------------------------------------------------------------------------
$newDate = "23/01/2017"
$ppPlaceholderDate = 16

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    $phType = -1
    try { $phType = $shp.PlaceholderFormat.Type } catch {}
    if ($phType -eq $ppPlaceholderDate) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $lyt = $layouts.Item($L)
    for ($i = 1; $i -le $lyt.Shapes.Count; $i++) {
        $shp = $lyt.Shapes.Item($i)
        $phType = -1
        try { $phType = $shp.PlaceholderFormat.Type } catch {}
        if ($phType -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Nudge the "ZoneTexte 25" ("File d'attente") label on slide 8 to its new
#    position: off x=6655866,y=2243274 (EMU) -> x=4791044,y=2065091 (EMU).
#    Shape.Left/.Top are expressed in points (1 pt = 12700 EMU).
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
for ($i = 1; $i -le $s8.Shapes.Count; $i++) {
    $shp = $s8.Shapes.Item($i)
    if ($shp.Name -eq "ZoneTexte 25") {
        $shp.Left = 377.2476
        $shp.Top = 162.6056
        break
    }
}
